$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H51").Value = 17537.857
$ws.Range("J51").Value = 12479.25
$ws.Range("L51").Value = 12479.25
$ws.Range("N51").Value = -13447.25
$ws.Range("H103").Value = 434.64706
$ws.Range("I103").Value = 269.5263
$ws.Range("J103").Value = 917.3077
$ws.Range("K103").Value = 808.5789
$ws.Range("L103").Value = 2751.9231
$ws.Range("M103").Value = -222.5789
$ws.Range("N103").Value = -3923.9231
$ws.Range("H116").Value = 14712780
$ws.Range("I116").Value = 27781254
$ws.Range("J116").Value = 10747.625
$ws.Range("K116").Value = 27781254
$ws.Range("L116").Value = 10747.625
$ws.Range("M116").Value = -27777812
$ws.Range("N116").Value = -17631.625
$ws.Range("H132").Value = 3702.1765
$ws.Range("I132").Value = 2912.25
$ws.Range("K132").Value = 8736.75
$ws.Range("M132").Value = -6206.75

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 3101.81
$ws.Range("I32").Value = 2954.922
$ws.Range("K32").Value = 2954.922
$ws.Range("M32").Value = -2667.922
$ws.Range("H61").Value = 5237.9517
$ws.Range("I61").Value = 2672.7874
$ws.Range("J61").Value = 13275.467
$ws.Range("K61").Value = 2672.7874
$ws.Range("L61").Value = 13275.467
$ws.Range("M61").Value = -2460.7874
$ws.Range("N61").Value = -13699.467
$ws.Range("H136").Value = 5237.9517
$ws.Range("I136").Value = 2672.7874
$ws.Range("J136").Value = 13275.467
$ws.Range("K136").Value = 8018.3622
$ws.Range("L136").Value = 39826.401
$ws.Range("M136").Value = -5468.3622
$ws.Range("N136").Value = -44926.401

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H80").Value = 33333832
$ws.Range("I80").Value = 83334136
$ws.Range("J80").Value = 295.77777
$ws.Range("K80").Value = 83334136
$ws.Range("L80").Value = 295.77777
$ws.Range("M80").Value = -83333138
$ws.Range("N80").Value = -2291.77777
$ws.Range("H83").Value = 33333832
$ws.Range("I83").Value = 83334136
$ws.Range("J83").Value = 295.77777
$ws.Range("K83").Value = 416670680
$ws.Range("L83").Value = 1478.88885
$ws.Range("M83").Value = -416665688
$ws.Range("N83").Value = -11462.88885
$ws.Range("H105").Value = 3638.8
$ws.Range("I105").Value = 3000
$ws.Range("J105").Value = 3798.5
$ws.Range("K105").Value = 3000
$ws.Range("L105").Value = 3798.5
$ws.Range("M105").Value = -1253
$ws.Range("N105").Value = -7292.5
$ws.Range("H107").Value = 62506652
$ws.Range("I107").Value = 86543900
$ws.Range("K107").Value = 86543900
$ws.Range("M107").Value = -86541980

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 8228.857
$ws.Range("I31").Value = 3676.2856
$ws.Range("J31").Value = 11643.286
$ws.Range("K31").Value = 3676.2856
$ws.Range("L31").Value = 11643.286
$ws.Range("M31").Value = -3381.2856
$ws.Range("N31").Value = -12233.286
$ws.Range("H34").Value = 8228.857
$ws.Range("I34").Value = 3676.2856
$ws.Range("J34").Value = 11643.286
$ws.Range("K34").Value = 3676.2856
$ws.Range("L34").Value = 11643.286
$ws.Range("M34").Value = -3474.2856
$ws.Range("N34").Value = -12047.286
$ws.Range("H58").Value = 16674336
$ws.Range("I58").Value = 50000828
$ws.Range("K58").Value = 50000828
$ws.Range("M58").Value = -50000625
$ws.Range("H99").Value = 8351.200000000001
$ws.Range("I99").Value = 10626.25
$ws.Range("J99").Value = 5751.143
$ws.Range("K99").Value = 10626.25
$ws.Range("L99").Value = 5751.143
$ws.Range("M99").Value = -9128.25
$ws.Range("N99").Value = -8747.143
$ws.Range("H126").Value = 8351.200000000001
$ws.Range("I126").Value = 10626.25
$ws.Range("J126").Value = 5751.143
$ws.Range("K126").Value = 31878.75
$ws.Range("L126").Value = 17253.429
$ws.Range("M126").Value = -29408.75
$ws.Range("N126").Value = -22193.429
$ws.Range("H136").Value = 16674336
$ws.Range("I136").Value = 50000828
$ws.Range("K136").Value = 150002484
$ws.Range("M136").Value = -149999934

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H107").Value = 25000246
$ws.Range("J107").Value = 33333530
$ws.Range("L107").Value = 100000590
$ws.Range("N107").Value = -100004430

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H135").Value = 89999
$ws.Range("J135").Value = 89999
$ws.Range("L135").Value = 89999
$ws.Range("N135").Value = -100139

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 8054.9165
$ws.Range("I7").Value = 5883
$ws.Range("J7").Value = 8778.888999999999
$ws.Range("K7").Value = 5883
$ws.Range("L7").Value = 8778.888999999999
$ws.Range("M7").Value = -5771
$ws.Range("N7").Value = -9002.888999999999
$ws.Range("H22").Value = 2959
$ws.Range("I22").Value = 1460
$ws.Range("K22").Value = 1460
$ws.Range("M22").Value = -1165
$ws.Range("H27").Value = 2959
$ws.Range("I27").Value = 1460
$ws.Range("K27").Value = 1460
$ws.Range("M27").Value = -1353
$ws.Range("H46").Value = 4111
$ws.Range("J46").Value = 4111
$ws.Range("L46").Value = 4111
$ws.Range("N46").Value = -4487
$ws.Range("H55").Value = 50000510
$ws.Range("I55").Value = 200000110
$ws.Range("J55").Value = 644.4666999999999
$ws.Range("K55").Value = 200000110
$ws.Range("L55").Value = 644.4666999999999
$ws.Range("M55").Value = -199999937
$ws.Range("N55").Value = -990.4666999999999
$ws.Range("H122").Value = 5015.8696
$ws.Range("I122").Value = 3335
$ws.Range("J122").Value = 8857.857
$ws.Range("K122").Value = 10005
$ws.Range("L122").Value = 26573.571
$ws.Range("M122").Value = -7555
$ws.Range("N122").Value = -31473.571
$ws.Range("H126").Value = 8054.9165
$ws.Range("I126").Value = 5883
$ws.Range("J126").Value = 8778.888999999999
$ws.Range("K126").Value = 17649
$ws.Range("L126").Value = 26336.667
$ws.Range("M126").Value = -15179
$ws.Range("N126").Value = -31276.667

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H45").Value = 10936.857
$ws.Range("I45").Value = 6973.5
$ws.Range("J45").Value = 12522.2
$ws.Range("K45").Value = 6973.5
$ws.Range("L45").Value = 12522.2
$ws.Range("M45").Value = -6482.5
$ws.Range("N45").Value = -13504.2
$ws.Range("H122").Value = 6227.6113
$ws.Range("I122").Value = 5032.8335
$ws.Range("J122").Value = 6825
$ws.Range("K122").Value = 15098.5005
$ws.Range("L122").Value = 20475
$ws.Range("M122").Value = -12648.5005
$ws.Range("N122").Value = -25375
$ws.Range("H136").Value = 34488420
$ws.Range("J136").Value = 9425.3125
$ws.Range("L136").Value = 28275.9375
$ws.Range("N136").Value = -33375.9375
